# Update the cryptocurrency price table (rows 2-51, columns B/C/D/E)
# Column D holds price strings that look numeric (e.g. "1.000", "0.4826");
# force them to remain text so Excel doesn't reinterpret them as numbers,
# matching the original inlineStr/text storage in the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.466.87'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '1.869.60'
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("D5").Value = '236.41'
$ws.Range("E5").Value = '  -0.50%  '
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").Value = '0.4826'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '0.2804'
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.06513'
$ws.Range("E9").Value = '  +0.21%  '
$ws.Range("D10").Value = '1.889.33'
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("D11").Value = '0.07450'
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").Value = '16.27'
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").Value = '5.088'
$ws.Range("E13").Value = '  +0.27%  '
$ws.Range("D14").Value = '87.21'
$ws.Range("E14").Value = '  -0.46%  '
$ws.Range("D15").Value = '0.6426'
$ws.Range("E15").Value = '  -1.70%  '
$ws.Range("D16").Value = '30.444.98'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").Value = '13.02'
$ws.Range("E18").Value = '  -1.60%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '231.32'
$ws.Range("E19").Value = '  +4.98%  '
$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").Value = '0.000007493'
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.099.72'
$ws.Range("E21").Value = '  -0.64%  '
$ws.Range("B22").Value = 'BinanceUSD'
$ws.Range("C22").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.151'
$ws.Range("E23").Value = '  -2.60%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '6.107'
$ws.Range("E24").Value = '  -0.87%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '170.49'
$ws.Range("E25").Value = '  +1.89%  '
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").Value = '9.351'
$ws.Range("E26").Value = '  +0.54%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '18.35'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.906'
$ws.Range("E28").Value = '  -3.00%  '
$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").Value = '0.1051'
$ws.Range("E29").Value = '  +12.67%  '
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").Value = '1.385'
$ws.Range("E30").Value = '  -4.85%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '4.273'
$ws.Range("E31").Value = '  -0.62%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.990'
$ws.Range("E32").Value = '  -0.68%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.04984'
$ws.Range("E33").Value = '  -1.14%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.181'
$ws.Range("E34").Value = '  -1.90%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.7425'
$ws.Range("E35").Value = '  -0.41%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").Value = '0.9994'
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").Value = '2.710'
$ws.Range("E37").Value = '  -0.32%  '
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").Value = '0.01932'
$ws.Range("E38").Value = '  +5.69%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.634'
$ws.Range("E39").Value = '  +0.92%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '0.9172'
$ws.Range("E40").Value = '  +0.28%  '
$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = '2.049'
$ws.Range("E41").Value = '  -1.59%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '106.06'
$ws.Range("E42").Value = '  -0.70%  '
$ws.Range("B43").Value = 'PaxDollar'
$ws.Range("C43").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D43").Value = '0.9970'
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4198'
$ws.Range("E44").Value = '  -1.44%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '5.584'
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").Value = '7.226'
$ws.Range("E46").Value = '  -2.17%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").Value = '62.13'
$ws.Range("E47").Value = '  -2.52%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D48").Value = '0.1228'
$ws.Range("E48").Value = '  -4.83%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '8.895'
$ws.Range("E49").Value = '  -0.98%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '1.424'
$ws.Range("E50").Value = '  -3.35%  '
$ws.Range("B51").Value = 'Elrond'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D51").Value = '33.61'
$ws.Range("E51").Value = '  -0.43%  '

$priceRange.Style = "Normal"
